# lab3 (2).xlsx — fix the u(U) combined-uncertainty formulas for the U
# measurement range (column O / Q), per "OKURWA WSZYSTKIE U_B SA ZLE XDDDDDDDDD".
#
# The old model computed the resolution/device uncertainty ua(U)-style term in
# column O as a plain linear expression, and then divided it by SQRT(3) into
# column Q (rectangular-distribution divisor) to get u(U)[V].
#
# The corrected model folds both steps into a single RSS-style expression in
# column O directly: SQRT(POWER(0.05/100*D<row-49>+0.001,2)/(1.73)) — and the
# now-redundant column Q is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the corrected u(U)[V] column, placed right above where the
# corrected formula range starts (row 100).
$ws.Range("O99").Value = "u(U)[V]"

# Rewrite O100:O160 with the corrected combined-uncertainty formula; D<row-49>
# is the matching "Temp/100" style denominator used by the old formula too.
for ($r = 100; $r -le 160; $r++) {
    $dRow = $r - 49
    $ws.Cells.Item($r, 15).Formula = "=SQRT(POWER(0.05/100*D$dRow+0.001,2)/(1.73))"
}

# The old Q100:Q160 helper column (O/SQRT(3)) is obsolete now, drop it.
$ws.Range("Q100:Q160").ClearContents()

# Update the view: scroll so column D is leftmost, and leave the selection on
# the newly-edited O100 cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 4
$ws.Range("O100").Select()
